# Update input files to 2025 FE data
# The "Further education and skills achievements" ILR-sourced rows (5 and 6)
# get refreshed links (new dataset/permalink GUIDs) and refreshed
# latest/next period labels, rolled forward by one year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New hyperlinks first (matches the order the refreshed ILR links/periods
# were authored in), then the refreshed latest/next period labels.

# Row 5: Further education and skills achievements and participation by
# provision, level and age group
$ws.Range("B5").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-catalogue/data-set/1977cdbc-7eae-4257-a8c9-3281bb2dbfa9'>Individualised Learner Record</a>"

# Row 6: Further education and skills achievements by sector subject area
$ws.Range("B6").Value = "<a href='https://explore-education-statistics.service.gov.uk/data-tables/permalink/529ad7b1-7a0f-419c-eb81-08de29d3af56'>Individualised Learner Record</a>"

$ws.Range("C5").Value = "Aug 2024 – Jul 2025 (27/11/25)"
$ws.Range("D5").Value = "Aug 2025 – Jul 2026 (Nov 26)"
$ws.Range("C6").Value = "Aug 2024 – Jul 2025 (27/11/25)"
$ws.Range("D6").Value = "Aug 2025 – Jul 2026 (Nov 26)"

# Reflect the author's final on-screen selection (was B22, now A6) and drop
# the scrolled-down top-left cell left over from the previous save.
$ws.Range("A6").Select()
